# Adjust Standing Right Kick animation and frame data, allow simultaneous hits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Link the top "srk" computed table (rows 3-18) to pull its B column
#     (startup) from the mirrored input table below (rows 23-38), so both
#     tables stay in sync (row N <-> row N+20).
$ws.Range("B3").Formula = "=B23"
$ws.Range("B4:B18").Formula = "=B24"

# --- Frame-data tweak for Standing Right Kick (row 6 / row 26): startup
#     goes from 11 to 14, which shifts the onhit/onblock advantage values.
$ws.Range("B26").Value = 14
$ws.Range("F6").Value = 37
$ws.Range("G6").Value = 20

# --- Move the active selection/view (no more frozen scroll to A7; land on
#     D23 instead of F21).
[void]$ws.Range("D23").Select()
